$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q (17) formatting into new column R (18) for rows 3-34,
# mirroring the per-row style used in column Q.
for ($r = 3; $r -le 34; $r++) {
    $src = $ws.Cells.Item($r, 17)   # column Q
    $dst = $ws.Cells.Item($r, 18)   # column R
    $src.Copy($dst)
}

# Header cell R4 = 2021
$ws.Cells.Item(4, 18).Value = 2021

# Numeric data values for column R (2021 figures)
$ws.Cells.Item(5, 18).Value = 109
$ws.Cells.Item(6, 18).Value = 74
$ws.Cells.Item(7, 18).Value = 35
$ws.Cells.Item(8, 18).Value = 36
$ws.Cells.Item(9, 18).Value = 35
$ws.Cells.Item(10, 18).Value = 1
$ws.Cells.Item(11, 18).Value = 15
$ws.Cells.Item(12, 18).Value = 8
$ws.Cells.Item(13, 18).Value = 7
$ws.Cells.Item(14, 18).Value = 12
$ws.Cells.Item(15, 18).Value = 7
$ws.Cells.Item(16, 18).Value = 5
$ws.Cells.Item(17, 18).Value = "-"
$ws.Cells.Item(18, 18).Value = "-"
$ws.Cells.Item(19, 18).Value = "-"
$ws.Cells.Item(20, 18).Value = 17
$ws.Cells.Item(21, 18).Value = 8
$ws.Cells.Item(22, 18).Value = 9
$ws.Cells.Item(23, 18).Value = 9
$ws.Cells.Item(24, 18).Value = 7
$ws.Cells.Item(25, 18).Value = 2
$ws.Cells.Item(26, 18).Value = 20
$ws.Cells.Item(27, 18).Value = 9
$ws.Cells.Item(28, 18).Value = 11
$ws.Cells.Item(29, 18).Value = "-"
$ws.Cells.Item(30, 18).Value = "-"
$ws.Cells.Item(31, 18).Value = "-"
$ws.Cells.Item(32, 18).Value = "-"
$ws.Cells.Item(33, 18).Value = "-"
$ws.Cells.Item(34, 18).Value = "-"
# R3 stays empty (only formatting copied, matching <c r="R3" s="8"/>)

# Update selection to match the post-edit state (activeCell R1)
$ws.Range("R1").Select()

Write-Output "done"
